$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "db" (sheet1) -- reorganize the lookup table:
#  * merge SNU_ESTATAL / SNU_PRIVADO into a single SNU row
#  * rename EST_SERV -> EST_SERVICIO
#  * split "Espacios religiosos" into CULTO_CATOLICO / CULTO_OTRO
#  * drop the stray column H helper cells (H11, H21)
#  * move the yellow highlight to the new row positions
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("db")

# Clear the two orphan helper cells that lived in column H; this also
# shrinks the sheet's used range from A1:H25 down to A1:E25.
$ws.Range("H11").ClearContents()
$ws.Range("H21").ClearContents()

# Row 2-7 keep their values (Inicial/Primario/Secundario Estatal/Privado).
# Row 8 now holds the merged SNU entry instead of SNU_ESTATAL.
$ws.Range("A8").Value = "SNU"
$ws.Range("B8").Value = "Nivel SNU (Estatal y Privado)"
$ws.Range("C8").Value = "EDU"
$ws.Range("D8").Value = "Educación"
$ws.Range("E8").Value = 1

# Row 9 now holds what used to be row 10 (Universidad).
$ws.Range("A9").Value = "UNI"
$ws.Range("B9").Value = "Universidad"
$ws.Range("C9").Value = "EDU"
$ws.Range("D9").Value = "Educación"
$ws.Range("E9").Value = 1

# Row 10 now holds what used to be row 11 (Salud).
$ws.Range("A10").Value = "EDIFICIO_DE_SALUD    "
$ws.Range("B10").Value = "Institución de Atención Médica de Gestión Estatal"
$ws.Range("C10").Value = "SAL"
$ws.Range("D10").Value = "Salud y Cuidado"
$ws.Range("E10").Value = 1

# Row 11: EST_SERV renamed to EST_SERVICIO (label unchanged).
$ws.Range("A11").Value = "EST_SERVICIO"
$ws.Range("B11").Value = "Estación de servicio"
$ws.Range("C11").Value = "TRA"
$ws.Range("D11").Value = "Transporte y espacio público"
$ws.Range("E11").Value = 1

# Row 12 now holds what used to be row 13 (Autopista).
$ws.Range("A12").Value = "AUTOPISTA"
$ws.Range("B12").Value = "Acceso a Autopista o vía principal"
$ws.Range("C12").Value = "TRA"
$ws.Range("D12").Value = "Transporte y espacio público"
$ws.Range("E12").Value = 2

# Row 13 now holds what used to be row 14 (Colectivo Municipal).
$ws.Range("A13").Value = "COLEC_MUNICIPAL"
$ws.Range("B13").Value = "Parada de colectivo de jurisdicción Municipal"
$ws.Range("C13").Value = "TRA"
$ws.Range("D13").Value = "Transporte y espacio público"
$ws.Range("E13").Value = 2

# Row 14 now holds what used to be row 15 (Colectivo Provincial).
$ws.Range("A14").Value = "COLEC_PROVINCIAL"
$ws.Range("B14").Value = "Parada de colectivo de jurisdicción Provincial"
$ws.Range("C14").Value = "TRA"
$ws.Range("D14").Value = "Transporte y espacio público"
$ws.Range("E14").Value = 2

# Row 15 now holds what used to be row 16 (Colectivo Nacional).
$ws.Range("A15").Value = "COLEC_NACIONAL"
$ws.Range("B15").Value = "Parada de colectivo de jurisdicción Nacional"
$ws.Range("C15").Value = "TRA"
$ws.Range("D15").Value = "Transporte y espacio público"
$ws.Range("E15").Value = 2

# Row 16 now holds what used to be row 17 (Tren/Subte).
$ws.Range("A16").Value = "TREN_EST"
$ws.Range("B16").Value = "Estación de Tren/Subte/Premetro"
$ws.Range("C16").Value = "TRA"
$ws.Range("D16").Value = "Transporte y espacio público"
$ws.Range("E16").Value = 1

# Row 17 now holds what used to be row 18 (Espacio verde) and keeps the
# highlight, now additionally applied to column A as well.
$ws.Range("A17").Value = "ESPACIO_VERDE "
$ws.Range("B17").Value = "Parques y Espacio verdes"
$ws.Range("C17").Value = "OFP"
$ws.Range("D17").Value = "Atención estatal, Seguridad y Justicia"
$ws.Range("E17").Value = 3

# Row 18: first half of the old "Espacios religiosos" split (Católico).
$ws.Range("A18").Value = "CULTO_CATOLICO"
$ws.Range("B18").Value = "Centro Religioso Católico"
$ws.Range("C18").Value = "OTR"
$ws.Range("D18").Value = "Redes de sociabilidad"
$ws.Range("E18").Value = 1

# Row 19: second half of the split (No-Católico / Evangelista).
$ws.Range("A19").Value = "CULTO_OTRO"
$ws.Range("B19").Value = "Centro Religioso No-Católico (Evangelista)"
$ws.Range("C19").Value = "OTR"
$ws.Range("D19").Value = "Redes de sociabilidad"
$ws.Range("E19").Value = 1

# Rows 20-25 are unchanged in content (Police, Bank, Super, Fabrica,
# Planta Transformadora, Cel) -- re-assert them defensively so the sheet
# is fully deterministic regardless of the starting state.
$ws.Range("A20").Value = "POLICE"
$ws.Range("B20").Value = "Comisaria"
$ws.Range("C20").Value = "OFP"
$ws.Range("D20").Value = "Atención estatal, Seguridad y Justicia"
$ws.Range("E20").Value = 1

$ws.Range("A21").Value = "BANK"
$ws.Range("B21").Value = "Bancos"
$ws.Range("C21").Value = "PRO"
$ws.Range("D21").Value = "Producción, comercio y sistema financiero"
$ws.Range("E21").Value = 1

$ws.Range("A22").Value = "SUPER"
$ws.Range("B22").Value = "Supermercado"
$ws.Range("C22").Value = "PRO"
$ws.Range("D22").Value = "Producción, comercio y sistema financiero"
$ws.Range("E22").Value = 1

$ws.Range("A23").Value = "FABRICA              "
$ws.Range("B23").Value = "Complejo fabril"
$ws.Range("C23").Value = "PRO"
$ws.Range("D23").Value = "Producción, comercio y sistema financiero"
$ws.Range("E23").Value = 1

$ws.Range("A24").Value = "PLANTA_TRANSFORMADORA"
$ws.Range("B24").Value = "Transformador eléctrico de media tensión"
$ws.Range("C24").Value = "INF"
$ws.Range("D24").Value = "Infraestructura"
$ws.Range("E24").Value = 1

$ws.Range("A25").Value = "CEL"
$ws.Range("B25").Value = "Antena de celular"
$ws.Range("C25").Value = "INF"
$ws.Range("D25").Value = "Infraestructura"
$ws.Range("E25").Value = 1

# Re-paint the yellow highlight so it lands exactly on the rows that hold
# it after the re-shuffle: A12, E12:E15, A17, E17.
$highlightRanges = @("E13","E14","E15","E16","E18")
foreach ($addr in $highlightRanges) {
    $ws.Range($addr).Interior.Pattern = -4142
}
$ws.Range("A12").Interior.Color = 65535
$ws.Range("E12").Interior.Color = 65535
$ws.Range("E13").Interior.Color = 65535
$ws.Range("E14").Interior.Color = 65535
$ws.Range("E15").Interior.Color = 65535
$ws.Range("A17").Interior.Color = 65535
$ws.Range("E17").Interior.Color = 65535

# Move the cursor/selection and update the view like the authored file.
$ws.Range("B29").Select()

# Add the new reviewer comment on A10 (Farmacia row) describing the
# request to add private pharmacies and hospitals.
$comment = $ws.Range("A10").AddComment("Autor:`nAgregar farmacias y Sanatorios privados.")

# ---------------------------------------------------------------------------
# Sheet "DRAFT" (sheet2) -- no data changed, only move the cursor.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DRAFT")
$ws2.Range("E8").Select()
